# Coffee Maker Acceptance Tests.docx edit
#
# 1) "Precondition: addInventory1 has run successfully." ->
#       "Precondition:" + " Run CoffeeMaker"   (split into two runs)
# 2) "Coffee: 16, Milk: 17, Sugar: 16, Chocolate: 17" ->
#       "Coffee: 15, Milk: 15, Sugar: 15" + ", Chocolate: 1" + "5"
#       (split into three runs)
#
# Word normally coalesces text assigned to adjacent Ranges that share
# identical run formatting back into a single <w:r>. Toggling a
# character-formatting property (Bold on, then back off) on the
# boundary forces the engine to keep the runs distinct, which is what
# lets us reproduce the multi-run XML shape from the diff.

$d = $word.ActiveDocument

# --- Change 1: Precondition line -------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Precondition: addInventory1 has run successfully.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$s1 = $rng1.Start
$e1 = $rng1.End

$firstLen1 = "Precondition:".Length

# Replace the trailing portion of the sentence with the new wording.
$tail1 = $d.Range($s1 + $firstLen1, $e1)
$tail1.Text = " Run CoffeeMaker"

# Force a run boundary right after "Precondition:" so the new text
# lands in its own <w:r> instead of merging with the first run.
$newTailEnd1 = $s1 + $firstLen1 + " Run CoffeeMaker".Length
$split1 = $d.Range($s1 + $firstLen1, $newTailEnd1)
$split1.Bold = 1
$split1.Bold = 0

# --- Change 2: inventory counts line ----------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Coffee: 16, Milk: 17, Sugar: 16, Chocolate: 17", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$s2 = $rng2.Start
$e2 = $rng2.End

$whole2 = $d.Range($s2, $e2)
$whole2.Text = "Coffee: 15, Milk: 15, Sugar: 15, Chocolate: 15"

$part1Len2 = "Coffee: 15, Milk: 15, Sugar: 15".Length
$part2Len2 = ", Chocolate: 1".Length

# Boundary between part 1 and part 2.
$splitA2 = $d.Range($s2 + $part1Len2, $s2 + $part1Len2 + $part2Len2)
$splitA2.Bold = 1
$splitA2.Bold = 0

# Boundary between part 2 and part 3.
$splitB2 = $d.Range($s2 + $part1Len2 + $part2Len2, $s2 + $part1Len2 + $part2Len2 + 1)
$splitB2.Bold = 1
$splitB2.Bold = 0
